$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43/44: Hedera and dogwifhat swap positions
$ws.Cells.Item(43, "B").Value = "dogwifhat"
$ws.Cells.Item(43, "C").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(43, "D").Value = "'2.33"
$ws.Cells.Item(43, "E").Value = "  -0.06%  "
$ws.Cells.Item(44, "B").Value = "Hedera"
$ws.Cells.Item(44, "C").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(44, "D").Value = "'0.0619"
$ws.Cells.Item(44, "E").Value = "  +1.46%  "

# Updated prices and volume percentages
$ws.Cells.Item(2, "D").Value = "65.882.71"
$ws.Cells.Item(2, "E").Value = "  +0.31%  "
$ws.Cells.Item(3, "D").Value = "2.667.27"
$ws.Cells.Item(3, "E").Value = "  -0.45%  "
$ws.Cells.Item(4, "E").Value = "  +0.01%  "
$ws.Cells.Item(5, "D").Value = "'599.31"
$ws.Cells.Item(5, "E").Value = "  -0.28%  "
$ws.Cells.Item(6, "D").Value = "'157.96"
$ws.Cells.Item(6, "E").Value = "  +0.54%  "
$ws.Cells.Item(7, "D").Value = "'0.651"
$ws.Cells.Item(8, "E").Value = "  +0.02%  "
$ws.Cells.Item(9, "D").Value = "'0.127"
$ws.Cells.Item(9, "E").Value = "  -2.47%  "
$ws.Cells.Item(10, "D").Value = "'0.403"
$ws.Cells.Item(10, "E").Value = "  +0.18%  "
$ws.Cells.Item(11, "D").Value = "'5.86"
$ws.Cells.Item(11, "E").Value = "  -0.05%  "
$ws.Cells.Item(12, "E").Value = "  +1.58%  "
$ws.Cells.Item(13, "D").Value = "'29.12"
$ws.Cells.Item(13, "E").Value = "  -0.94%  "
$ws.Cells.Item(14, "E").Value = "  -2.47%  "
$ws.Cells.Item(15, "D").Value = "3.143.69"
$ws.Cells.Item(15, "E").Value = "  -0.45%  "
$ws.Cells.Item(16, "D").Value = "65.770.99"
$ws.Cells.Item(16, "E").Value = "  +0.39%  "
$ws.Cells.Item(17, "D").Value = "2.627.53"
$ws.Cells.Item(17, "E").Value = "  -1.66%  "
$ws.Cells.Item(18, "E").Value = "  -1.18%  "
$ws.Cells.Item(19, "D").Value = "'4.80"
$ws.Cells.Item(19, "E").Value = "  +0.01%  "
$ws.Cells.Item(20, "D").Value = "'351.57"
$ws.Cells.Item(20, "E").Value = "  -0.11%  "
$ws.Cells.Item(21, "E").Value = "  -1.47%  "
$ws.Cells.Item(22, "E").Value = "  -0.09%  "
$ws.Cells.Item(23, "E").Value = "  +0.35%  "
$ws.Cells.Item(24, "D").Value = "'1.84"
$ws.Cells.Item(24, "E").Value = "  +11.65%  "
$ws.Cells.Item(25, "D").Value = "'0.0000113"
$ws.Cells.Item(25, "E").Value = "  +0.29%  "
$ws.Cells.Item(26, "D").Value = "'9.65"
$ws.Cells.Item(26, "E").Value = "  -0.37%  "
$ws.Cells.Item(27, "E").Value = "  +2.12%  "
$ws.Cells.Item(28, "D").Value = "'567.29"
$ws.Cells.Item(28, "E").Value = "  +6.00%  "
$ws.Cells.Item(29, "D").Value = "'8.24"
$ws.Cells.Item(29, "E").Value = "  +1.54%  "
$ws.Cells.Item(30, "E").Value = "  -2.12%  "
$ws.Cells.Item(31, "D").Value = "'0.999"
$ws.Cells.Item(31, "E").Value = "  -0.13%  "
$ws.Cells.Item(32, "D").Value = "'2.15"
$ws.Cells.Item(32, "E").Value = "  -0.12%  "
$ws.Cells.Item(33, "E").Value = "  +3.89%  "
$ws.Cells.Item(34, "D").Value = "'6.69"
$ws.Cells.Item(34, "E").Value = "  +3.38%  "
$ws.Cells.Item(35, "E").Value = "  +0.81%  "
$ws.Cells.Item(36, "E").Value = "  -0.32%  "
$ws.Cells.Item(37, "D").Value = "'20.58"
$ws.Cells.Item(37, "E").Value = "  +0.36%  "
$ws.Cells.Item(38, "E").Value = "  -0.11%  "
$ws.Cells.Item(39, "E").Value = "  +0.64%  "
$ws.Cells.Item(40, "D").Value = "'154.05"
$ws.Cells.Item(40, "E").Value = "  -2.61%  "
$ws.Cells.Item(41, "D").Value = "'161.24"
$ws.Cells.Item(41, "E").Value = "  -2.12%  "
$ws.Cells.Item(42, "E").Value = "  -1.36%  "
$ws.Cells.Item(45, "D").Value = "'23.00"
$ws.Cells.Item(45, "E").Value = "  +0.64%  "
$ws.Cells.Item(46, "E").Value = "  +0.73%  "
$ws.Cells.Item(47, "E").Value = "  -1.46%  "
$ws.Cells.Item(48, "E").Value = "  +0.92%  "
$ws.Cells.Item(49, "D").Value = "'19.86"
$ws.Cells.Item(49, "E").Value = "  -1.65%  "
$ws.Cells.Item(50, "E").Value = "  -5.55%  "
$ws.Cells.Item(51, "E").Value = "  -0.94%  "
